$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated pivot values (rows 2-15), regenerated to include credit card
# account refund/credit transactions.
$data = @{
    2  = @(-81801.64, -38293.06, -60928.02, -181022.72)
    3  = @(126416.96, 46505.97, 165690.46, 338613.39)
    4  = @(-66810.42999999999, -38819.41, -73376.16, -179006)
    5  = @(892750.23, 521104.54, 840546.9, 2254401.67)
    6  = @(-32211.18, -15438.74, -33069.4, -80719.32000000001)
    7  = @(-82469.08, -31167.05, -76872.91, -190509.04)
    8  = @(-32168.74, -16132.08, -31931.7, -80232.52)
    9  = @(-102411.74, -54028.1, -123307.01, -279746.85)
    10 = @(-58558.86, -27644.43, -54895.93, -141099.22)
    11 = @(384725.52, 185434.74, 395603.59, 965763.85)
    12 = @(-384725.52, -185434.74, -395603.59, -965763.85)
    13 = @(-34618.25, -15629.56, -32754.63, -83002.44)
    14 = @(157971.31, 76187.89999999999, 179165.33, 413324.54)
    15 = @(686088.58, 406645.98, 698266.9300000001, 1791001.49)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Range("B$row").Value = $vals[0]
    $ws.Range("C$row").Value = $vals[1]
    $ws.Range("D$row").Value = $vals[2]
    $ws.Range("E$row").Value = $vals[3]
}
